# Auto-applies numeric value updates captured in the source diff.
# Each worksheet's table contains static market-data snapshot values
# (currentAveragePrice*, LevePrice*, LeveProfit*) with no formulas,
# so the edit simply overwrites the affected cells with their new values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3700
$ws.Range("J74").Value = 4000
$ws.Range("L74").Value = 4000
$ws.Range("N74").Value = -5872
$ws.Range("H77").Value = 3700
$ws.Range("J77").Value = 4000
$ws.Range("L77").Value = 20000
$ws.Range("N77").Value = -29360
$ws.Range("H113").Value = 3125.8572
$ws.Range("J113").Value = 3218
$ws.Range("L113").Value = 3218
$ws.Range("N113").Value = -9726
$ws.Range("H132").Value = 5292590.5
$ws.Range("I132").Value = 6804378
$ws.Range("J132").Value = 1334.3334
$ws.Range("K132").Value = 20413134
$ws.Range("L132").Value = 4003.0002
$ws.Range("M132").Value = -20410604
$ws.Range("N132").Value = -9063.0002
$ws.Range("H137").Value = 1342.2858
$ws.Range("I137").Value = 1178.5714
$ws.Range("J137").Value = 1997.1428
$ws.Range("K137").Value = 3535.7142
$ws.Range("L137").Value = 5991.428400000001
$ws.Range("M137").Value = -985.7142000000003
$ws.Range("N137").Value = -11091.4284
$ws.Range("H138").Value = 2470.5
$ws.Range("I138").Value = 3240
$ws.Range("J138").Value = 2278.125
$ws.Range("K138").Value = 9720
$ws.Range("L138").Value = 6834.375
$ws.Range("M138").Value = -4580
$ws.Range("N138").Value = -17114.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1099.1562
$ws.Range("I2").Value = 845.2632
$ws.Range("J2").Value = 1470.2307
$ws.Range("K2").Value = 845.2632
$ws.Range("L2").Value = 1470.2307
$ws.Range("M2").Value = -732.2632
$ws.Range("N2").Value = -1696.2307
$ws.Range("H74").Value = 1953.6957
$ws.Range("I74").Value = 1624.1666
$ws.Range("K74").Value = 1624.1666
$ws.Range("M74").Value = -750.1666
$ws.Range("H77").Value = 1953.6957
$ws.Range("I77").Value = 1624.1666
$ws.Range("K77").Value = 8120.833000000001
$ws.Range("M77").Value = -3752.833000000001
$ws.Range("H116").Value = 1099.1562
$ws.Range("I116").Value = 845.2632
$ws.Range("J116").Value = 1470.2307
$ws.Range("K116").Value = 845.2632
$ws.Range("L116").Value = 1470.2307
$ws.Range("M116").Value = 1448.7368
$ws.Range("N116").Value = -6058.2307
$ws.Range("H132").Value = 5332.185
$ws.Range("I132").Value = 5816.227
$ws.Range("J132").Value = 3202.4
$ws.Range("K132").Value = 17448.681
$ws.Range("L132").Value = 9607.200000000001
$ws.Range("M132").Value = -14918.681
$ws.Range("N132").Value = -14667.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1099.1562
$ws.Range("I3").Value = 845.2632
$ws.Range("J3").Value = 1470.2307
$ws.Range("K3").Value = 845.2632
$ws.Range("L3").Value = 1470.2307
$ws.Range("M3").Value = -731.2632
$ws.Range("N3").Value = -1698.2307
$ws.Range("H86").Value = 1486.2693
$ws.Range("I86").Value = 1313
$ws.Range("J86").Value = 1956.5714
$ws.Range("K86").Value = 1313
$ws.Range("L86").Value = 1956.5714
$ws.Range("M86").Value = -190
$ws.Range("N86").Value = -4202.5714
$ws.Range("H89").Value = 1486.2693
$ws.Range("I89").Value = 1313
$ws.Range("J89").Value = 1956.5714
$ws.Range("K89").Value = 6565
$ws.Range("L89").Value = 9782.857
$ws.Range("M89").Value = -949
$ws.Range("N89").Value = -21014.857
$ws.Range("H107").Value = 1933.3334
$ws.Range("I107").Value = 1400
$ws.Range("K107").Value = 1400
$ws.Range("M107").Value = 520
$ws.Range("H134").Value = 102171.2
$ws.Range("I134").Value = 168118.67
$ws.Range("J134").Value = 3250
$ws.Range("K134").Value = 504356.01
$ws.Range("L134").Value = 9750
$ws.Range("M134").Value = -501821.01
$ws.Range("N134").Value = -14820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 822.5
$ws.Range("I16").Value = 766.5
$ws.Range("J16").Value = 906.5
$ws.Range("K16").Value = 766.5
$ws.Range("L16").Value = 906.5
$ws.Range("M16").Value = -479.5
$ws.Range("N16").Value = -1480.5
$ws.Range("H31").Value = 3705647
$ws.Range("I31").Value = 1913.4894
$ws.Range("J31").Value = 28573572
$ws.Range("K31").Value = 1913.4894
$ws.Range("L31").Value = 28573572
$ws.Range("M31").Value = -1618.4894
$ws.Range("N31").Value = -28574162
$ws.Range("H34").Value = 3705647
$ws.Range("I34").Value = 1913.4894
$ws.Range("J34").Value = 28573572
$ws.Range("K34").Value = 1913.4894
$ws.Range("L34").Value = 28573572
$ws.Range("M34").Value = -1711.4894
$ws.Range("N34").Value = -28573976
$ws.Range("H41").Value = 4485.4546
$ws.Range("J41").Value = 4485.4546
$ws.Range("L41").Value = 4485.4546
$ws.Range("N41").Value = -5341.4546
$ws.Range("H50").Value = 16860
$ws.Range("J50").Value = 16860
$ws.Range("L50").Value = 16860
$ws.Range("N50").Value = -18110
$ws.Range("H51").Value = 12040
$ws.Range("J51").Value = 12040
$ws.Range("L51").Value = 12040
$ws.Range("N51").Value = -13512
$ws.Range("H59").Value = 31260
$ws.Range("J59").Value = 31260
$ws.Range("L59").Value = 31260
$ws.Range("N59").Value = -33550
$ws.Range("H60").Value = 12560.92
$ws.Range("J60").Value = 12721.792
$ws.Range("L60").Value = 12721.792
$ws.Range("N60").Value = -13743.792
$ws.Range("H61").Value = 12040
$ws.Range("J61").Value = 12040
$ws.Range("L61").Value = 12040
$ws.Range("N61").Value = -12736
$ws.Range("H68").Value = 39075
$ws.Range("J68").Value = 39075
$ws.Range("L68").Value = 39075
$ws.Range("N68").Value = -40573
$ws.Range("H71").Value = 39075
$ws.Range("J71").Value = 39075
$ws.Range("L71").Value = 117225
$ws.Range("N71").Value = -124713
$ws.Range("H74").Value = 33800
$ws.Range("J74").Value = 33800
$ws.Range("L74").Value = 33800
$ws.Range("N74").Value = -35548
$ws.Range("H77").Value = 33800
$ws.Range("J77").Value = 33800
$ws.Range("L77").Value = 101400
$ws.Range("N77").Value = -110136
$ws.Range("H99").Value = 1968.56
$ws.Range("I99").Value = 1458.3334
$ws.Range("J99").Value = 2439.5386
$ws.Range("K99").Value = 1458.3334
$ws.Range("L99").Value = 2439.5386
$ws.Range("M99").Value = 39.66660000000002
$ws.Range("N99").Value = -5435.5386
$ws.Range("H113").Value = 822.5
$ws.Range("I113").Value = 766.5
$ws.Range("J113").Value = 906.5
$ws.Range("K113").Value = 766.5
$ws.Range("L113").Value = 906.5
$ws.Range("M113").Value = 1403.5
$ws.Range("N113").Value = -5246.5
$ws.Range("H126").Value = 1968.56
$ws.Range("I126").Value = 1458.3334
$ws.Range("J126").Value = 2439.5386
$ws.Range("K126").Value = 4375.0002
$ws.Range("L126").Value = 7318.6158
$ws.Range("M126").Value = -1905.0002
$ws.Range("N126").Value = -12258.6158
$ws.Range("H132").Value = 3378.842
$ws.Range("I132").Value = 3092.3076
$ws.Range("J132").Value = 3999.6667
$ws.Range("K132").Value = 9276.9228
$ws.Range("L132").Value = 11999.0001
$ws.Range("M132").Value = -6746.9228
$ws.Range("N132").Value = -17059.0001
$ws.Range("H134").Value = 865.72095
$ws.Range("I134").Value = 824
$ws.Range("J134").Value = 1182.8
$ws.Range("K134").Value = 2472
$ws.Range("L134").Value = 3548.4
$ws.Range("M134").Value = 63
$ws.Range("N134").Value = -8618.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 754.75
$ws.Range("I5").Value = 495.3846
$ws.Range("K5").Value = 1486.1538
$ws.Range("M5").Value = -1374.1538
$ws.Range("H122").Value = 795629.2
$ws.Range("I122").Value = 2940.6
$ws.Range("J122").Value = 1236011.8
$ws.Range("K122").Value = 26465.4
$ws.Range("L122").Value = 11124106.2
$ws.Range("M122").Value = -24015.4
$ws.Range("N122").Value = -11129006.2
$ws.Range("H135").Value = 754.75
$ws.Range("I135").Value = 495.3846
$ws.Range("K135").Value = 4458.4614
$ws.Range("M135").Value = -1923.4614

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 50002044
$ws.Range("J122").Value = 2772.7273
$ws.Range("L122").Value = 8318.1819
$ws.Range("N122").Value = -13218.1819
$ws.Range("H132").Value = 102930.35
$ws.Range("I132").Value = 156093.53
$ws.Range("K132").Value = 468280.59
$ws.Range("M132").Value = -465750.59

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3752.111
$ws.Range("I132").Value = 4236.4287
$ws.Range("J132").Value = 3230.5386
$ws.Range("K132").Value = 12709.2861
$ws.Range("L132").Value = 9691.6158
$ws.Range("M132").Value = -10179.2861
$ws.Range("N132").Value = -14751.6158
$ws.Range("H136").Value = 5991.273
$ws.Range("I136").Value = 7363
$ws.Range("K136").Value = 22089
$ws.Range("M136").Value = -19539

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3400.4
$ws.Range("I132").Value = 2334.6667
$ws.Range("K132").Value = 7004.000100000001
$ws.Range("M132").Value = -4474.000100000001
